$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) — update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 86
$ws1.Range("F3").Value = 815
$ws1.Range("F4").Value = 47
$ws1.Range("F6").Value = 124
$ws1.Range("F8").Value = 4699
$ws1.Range("F10").Value = 5085
$ws1.Range("F11").Value = 584
$ws1.Range("F12").Value = 1277

# Sheet "全部类型" (sheet4) — update "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 86
$ws4.Range("F3").Value = 815
$ws4.Range("F4").Value = 47
$ws4.Range("F6").Value = 124
$ws4.Range("F9").Value = 4699
$ws4.Range("F11").Value = 5085
$ws4.Range("F12").Value = 584
$ws4.Range("F13").Value = 1277

$wb.Save()
